$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of column M (row 2-6) into new column N, and set new values
$ws.Range("N2").Value = $null
$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 6333
$ws.Range("N5").Value = 82675
$ws.Range("N6").Value = 300853

$ws.Range("M2:M6").Copy()
$ws.Range("N2:N6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("N3").Value = 2022
$ws.Range("N4").Value = 6333
$ws.Range("N5").Value = 82675
$ws.Range("N6").Value = 300853

$ws.Range("N2").Select()
